$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131; this shifts existing rows 131-155 down to 132-156
$ws.Rows.Item(131).Insert()

# Copy the row that is now at 132 (the original row 131) into the new row 131
# so that the columns that are identical across the whole data set
# (A, B, C, E, F, G, H, I, R) and the date formatting are carried over correctly.
$ws.Range("A132:R132").Copy($ws.Range("A131:R131"))

# Now overwrite the columns that hold the new record's specific data
$ws.Range("D131").Value = 44474
$ws.Range("J131").Value = 180
$ws.Range("K131").Value = 4500
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = 4500
$ws.Range("N131").Value = "$/docena de atados (3 kilos)"
$ws.Range("O131").Value = "Región Metropolitana"
$ws.Range("P131").Value = 1500
$ws.Range("Q131").Value = 3
